$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row at position 10. This shifts the former rows
#    10-14 (A 27955-2020, A 42532-2020, A 42512-2020, A 27642-2020,
#    A 38342-2020) down to rows 11-15.
# ------------------------------------------------------------------
$ws.Rows.Item(10).Insert()

# ------------------------------------------------------------------
# 2) Populate the newly inserted row 10 with the updated data for
#    case "A 38342-2020" (previously at row 14, now refreshed with
#    new survey numbers).
# ------------------------------------------------------------------
$ws.Range("A10").Value = "A 38342-2020"
$ws.Range("B10").Value = 44060
$ws.Range("C10").Value = 45203
$ws.Range("D10").Value = "NORRBOTTENS LÄN"
$ws.Range("E10").Value = "KIRUNA"
$ws.Range("F10").Value = "Allmännings- och besparingsskogar"
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 5
$ws.Range("P10").Value = 1
$ws.Range("Q10").Value = 5
$ws.Range("R10").Value = "Ostticka`r`nGranticka`r`nRosenticka`r`nTallriska`r`nUllticka"

$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KIRUNA/artfynd/A 38342-2020.xlsx", "A 38342-2020")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KIRUNA/kartor/A 38342-2020.png", "A 38342-2020")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KIRUNA/klagomål/A 38342-2020.docx", "A 38342-2020")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KIRUNA/klagomålsmail/A 38342-2020.docx", "A 38342-2020")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KIRUNA/tillsyn/A 38342-2020.docx", "A 38342-2020")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KIRUNA/tillsynsmail/A 38342-2020.docx", "A 38342-2020")'

# ------------------------------------------------------------------
# 3) The old "A 38342-2020" row (now duplicated at row 15 after the
#    insert shifted it down) is removed - its content now lives in
#    the refreshed row 10 above.
# ------------------------------------------------------------------
$ws.Rows.Item(15).Delete()

# ------------------------------------------------------------------
# 4) Across every data row (2-100), bump the "Förändrad" date
#    (column C) from 2023-10-03 (45202) to 2023-10-04 (45203).
# ------------------------------------------------------------------
for ($r = 2; $r -le 100; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

"done"
